$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column N: "Implement" (same bold style as the other
# header cells in row 2, so copy M2's formatting across first).
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("N2").Value = "Implement"

# Fill column N with the same rotating-name pattern as columns J:M,
# shifted one further column to the right (N = M shifted by one row).
$ws.Range("N3").Formula = "=M4"
$ws.Range("N4").Formula = "=M5"
$ws.Range("N5").Formula = "=M6"
$ws.Range("N6").Formula = "=M7"
$ws.Range("N7").Formula = "=M8"
$ws.Range("N8").Formula = "=M9"
$ws.Range("N9").Formula = "=M10"
$ws.Range("N10").Formula = "=M11"
$ws.Range("N11").Formula = "=M12"
$ws.Range("N12").Formula = "=M13"
$ws.Range("N13").Formula = "=M14"
$ws.Range("N14").Formula = "=M3"

# Column width recorded for column M as part of this edit.
$ws.Columns("M").ColumnWidth = 22.67

# Move the active selection to A14 (matches the post-edit saved view).
$ws.Range("A14").Select() | Out-Null
